# Append " (Changed main)" to the end of the first paragraph
# ("This is a Microsoft word document."), as three separate runs:
#   " ("  /  "Changed main"  /  ")"
#
# A plain Range.InsertAfter() sequence gets coalesced back into the
# existing run on save (same formatting => merged), so each new chunk is
# grown in its own temporary paragraph (via InsertParagraphAfter) and the
# paragraph mark that separates it from the previous text is then deleted.
# Word keeps the two runs distinct once the paragraphs are stitched back
# together, producing the separate <w:r> elements the diff expects.

$d = $word.ActiveDocument

function Append-AsNewRun($insertPos, $text) {
    $anchor = $d.Range($insertPos, $insertPos)
    $anchor.InsertParagraphAfter()
    $newPara = $d.Range($insertPos, $insertPos).Paragraphs(1).Next()
    $newPara.Range.InsertAfter($text)
    $mark = $d.Range($insertPos, $insertPos + 1)
    $mark.Delete()
}

$target = $d.Paragraphs(1)
$pos = $target.Range.End - 1

Append-AsNewRun $pos " ("
$pos = $pos + 2

Append-AsNewRun $pos "Changed main"
$pos = $pos + 12

Append-AsNewRun $pos ")"
